# Recompute column H ("客単価" / per-customer spend) on the
# "ABC分析_客構成" sheet as H = H / E (previous value divided by the
# customer-count column E, i.e. turning a per-segment total into a
# true per-customer average). Rows whose result is unchanged (E = 1)
# or whose H is already 0/blank are left untouched, matching the
# source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "ABC分析_客構成"

$newValues = @{
    2  = 1020.225903614458
    3  = 1103.491124260355
    4  = 920.436507936508
    5  = 1019.389312977099
    6  = 940.9856630824373
    7  = 15.92592592592593
    9  = 557.156862745098
    10 = 1017.930029154519
    11 = 730.899470899471
    12 = 569.4545454545455
    13 = 226.948051948052
    14 = 980
    15 = 188.8245614035088
    16 = 329.5555555555555
    17 = 306.4285714285714
    18 = 1393.286713286713
    19 = 111.625
    20 = 163.2
    21 = 844
    23 = 689.375
    26 = 471.2121212121212
    27 = 277.6
    29 = 1900
    33 = 1723.333333333333
    34 = 815.5555555555555
    39 = 787.3809523809524
    40 = 2266.25
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $newValues[$row]
}
